$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "level" column F, matching formatting of the neighbouring columns ---
# Header cell: same style as the other header cells (copy format from A1)
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "level"

# Data cells: same style as column A's numeric data cells (copy format from A2:A4)
$ws.Range("A2").Copy()
$ws.Range("F2:F4").PasteSpecial(-4122)
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 3

$excel.CutCopyMode = $false

# --- Update existing data (className / password columns changed content) ---
$ws.Range("B2").Value = "1к"
$ws.Range("B3").Value = "5к"
$ws.Range("B4").Value = "10к"

$ws.Range("E2").Value = "1к24"
$ws.Range("E3").Value = "1к24"
$ws.Range("E4").Value = "1к24"

# --- Update selection to match the new active cell ---
[void]$ws.Range("E8").Select()
